# Add the new "Debug_Messages" setting row to the settings sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Debug_Messages"
$ws.Range("B3").Value = "no"

# Leave the newly-added cell selected, matching the saved view state.
$ws.Range("B3").Select()
